# Modify Discount type entity ("Kedvezmény típus" sheet):
# add a new "szorzo" (multiplier) attribute row to the table.

$wb = $excel.ActiveWorkbook

$dstName  = "Kedvezmény típus"
$dst      = $wb.Worksheets.Item($dstName)
$srcTop   = $wb.Worksheets.Item("Légitársaság")   # donor for the plain "first data row" look
$srcBot   = $wb.Worksheets.Item("Kedvezmény")     # donor for the thick-bottom "last data row" look

# New attribute row values
$dst.Range("A3").Value = "szorzo"
$dst.Range("B3").Value = "NUMBER(1, 2)"
$dst.Range("C3").Value = "A kedvezményhez tartozó szorzó"

# Re-apply formatting so the table keeps a consistent look now that row 3
# is the new last row and row 2 is back to being a plain (non-bottom) row.
$srcTop.Range("A2:C2").Copy()
$dst.Range("A2:C2").PasteSpecial(-4122)   # xlPasteFormats

$srcBot.Range("A7:C7").Copy()
$dst.Range("A3:C3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# The workbook was saved while the "Jogosultság" tab was the active one.
$ws1 = $wb.Worksheets.Item("Jogosultság")
$ws1.Activate()
